$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All four "Rejected" status rows (I2, I8, I10, I16) become "Approved",
# and the corresponding "ReasonToReject" (column J) values are cleared out,
# matching row 17 which was already Approved with no rejection reason.
$ws.Range("I2").Value = "Approved"
$ws.Range("J2").ClearContents()

$ws.Range("I8").Value = "Approved"
$ws.Range("J8").ClearContents()

$ws.Range("I10").Value = "Approved"
$ws.Range("J10").ClearContents()

$ws.Range("I16").Value = "Approved"
$ws.Range("J16").ClearContents()

# Row 17 stays "Approved" (value unchanged, but shared-string table shrinks).
$ws.Range("I17").Value = "Approved"

# Update the saved view/selection state.
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("H21").Select()
